$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11 (shifts existing rows 11.. down by one)
$ws.Rows.Item(11).Insert()

# Populate the new row 11 - write SWGOHGG (col C) before Name (col A) so the
# shared-string table gains "7alximik7" before "Alximik", matching source order
$ws.Range("C11").Value = "7alximik7"
$ws.Range("A11").Value = "Alximik"
$ws.Range("B11").Value = ":flag_ru:"
$ws.Range("D11").Value = 15

# Grow Table1 by one row so it keeps covering the data (A1:D37 -> A1:D38)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D38")) | Out-Null

# Append a fresh blank row at the bottom of the sheet (row 46)
$ws.Range("A46:D46").Value = ""

# Move the selection to match the saved cursor position
$ws.Range("B12").Select() | Out-Null
